$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.194.22'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '3.153.28'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.19'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.25'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -3.81%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.150.88'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.29'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.21'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').Value = '3.678.55'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '3.159.45'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '63.114.82'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.68'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.05'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.06'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -4.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.701'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.71'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.68'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.01'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.71'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.17'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.99'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.96'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '52.75'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.82'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -3.66%  '
$ws.Range('D38').Value = '0.0₃0699'
$ws.Range('E38').Value = '  -7.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0389'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '419.85'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -4.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.74'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -8.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.29'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = '2.932.37'
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.111'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -5.82%  '
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.48'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.25'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -8.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.33'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -0.44%  '
